$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values, computed/regenerated to use K instead of Strike#.
$newK = @{
    2  = 7
    3  = 12
    4  = 10
    5  = 6
    6  = 4
    7  = 10
    8  = 6
    9  = 6
    10 = 10
    11 = 6
    12 = 7
    13 = 10
    14 = 10
    15 = 9
    16 = 3
    17 = 6
    18 = 6
    19 = 2
    20 = 5
    21 = 7
    22 = 4
    23 = 11
    24 = 5
    25 = 4
    26 = 9
    27 = 4
    28 = 1
    29 = 7
    30 = 6
    31 = 5
    32 = 7
    33 = 6
    34 = 3
    35 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
